$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted cell above so the new Date cell (G9) gets the
# same style (s="1", numFmtId 22) rather than creating a brand-new style.
$ws.Range("G8").Copy()
$ws.Range("G9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New row of trade data appended at row 9.
$ws.Range("A9").Value = 9390.1
$ws.Range("B9").Value = 9316.5
$ws.Range("C9").Value = 282.89999999999998
$ws.Range("D9").Value = 285.14
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = 0.79
$ws.Range("G9").Value = 42609.487280092595
$ws.Range("H9").Value = $true
